# error handling on scrapper_Service
#
# The scraped "part number" / "amount" sample row is replaced by the real
# scraped product data: a (numeric-looking) part number in A1 and a long
# alphanumeric part/description code in A2. Both labels are now shown in a
# bold, word-wrapped Calibri cell, the old "amount" column is gone, and the
# selection/row heights follow the new two-line layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: new scraped part number (stored as text, like the rest of the
#     scraped codes, even though it looks numeric) ---------------------
$ws.Range("A1").Value = "'3273114"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true

# --- B1: the old "amount" header is gone; clear it out -----------------
$ws.Range("B1").ClearContents()

# --- A2: long scraped part/description code ----------------------------
$ws.Range("A2").Value = "6SL32105BE211UV0"
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").WrapText = $true

# --- B2: old "amount" value cell removed entirely -----------------------
$ws.Range("B2").ClearContents()

# --- Row heights for the new (wrapped) two-row layout --------------------
$ws.Rows.Item(1).RowHeight = 15.65
$ws.Rows.Item(2).RowHeight = 44

# Keep the sheet's trailing padding row and add the extra blank row just
# above it, matching the re-saved sheet's sheetData tail.
$ws.Rows.Item(1048575).RowHeight = 12.8
$ws.Rows.Item(1048576).RowHeight = 12.8

# Selection moves to A2 after the scrape updates the sheet.
$ws.Range("A2").Select()
